$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header renames
$ws.Range("D1").Value = "addlDays"
$ws.Range("F1").Value = "leave"
$ws.Range("G1").Value = "leave"
$ws.Range("H1").Value = "leave"

# Row 2 (Akram Khan) leave detail cells - reformatted text
$ws.Range("F2").Value = "leaveType: C;`nstart: 1 Oct 2019;`nend: 11 Oct 2019;`nspecialDays: -2;"
$ws.Range("G2").Value = "leaveType: W1;`nstart: 12 Dec 2019;`nend: 17 Dec 2019;`nspecialDays: 4;"
$ws.Range("H2").Value = "leaveType: W2;`nstart: 12 Jan 2020;`nend: 14 Jan 2020;`nspecialDays: 0;"

# Row 3 (Ilyas Hussain) leave detail cells - reformatted text
$ws.Range("F3").Value = "leaveType: C;`nstart: 10 Sep 2019;`nend: 23 Sep 2019;`nspecialDays: -2;"
$ws.Range("G3").Value = "leaveType: W1;`nstart: 24 Dec 2019;`nend: 30 Dec 2019;`nspecialDays: 4;"

# Row 4 (Rafi Ullah) leave detail cells - reformatted text, G4 cleared
$ws.Range("F4").Value = "leaveType: C;`nstart: 2 Nov 2019;`nend: 16 Nov 2019;`nspecialDays: 0;"
$ws.Range("G4").Value = ""

# Update the selected cell to F2
$ws.Range("F2").Select() | Out-Null
